$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.145.51'
$ws.Range('E2').Value = '  -3.19%  '
$ws.Range('D3').Value = '3.134.78'
$ws.Range('E3').Value = '  -2.23%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''611.03'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').Value = '''149.35'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.30%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '3.132.61'
$ws.Range('E8').Value = '  -2.24%  '
$ws.Range('D9').Value = '''0.533'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('E10').Value = '  -4.75%  '
$ws.Range('D11').Value = '''5.57'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.49%  '
$ws.Range('E12').Value = '  -5.12%  '
$ws.Range('E13').Value = '  -3.92%  '
$ws.Range('E14').Value = '  -4.32%  '
$ws.Range('D15').Value = '3.612.63'
$ws.Range('E15').Value = '  -3.11%  '
$ws.Range('D16').Value = '64.172.65'
$ws.Range('E16').Value = '  -3.33%  '
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').Value = '3.132.65'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('E19').Value = '  -4.66%  '
$ws.Range('D20').Value = '''484.61'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.78%  '
$ws.Range('D21').Value = '''14.67'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.31%  '
$ws.Range('E22').Value = '  -2.93%  '
$ws.Range('E23').Value = '  -2.98%  '
$ws.Range('D24').Value = '''13.83'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -5.60%  '
$ws.Range('D25').Value = '''84.38'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '''2.95'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.96%  '
$ws.Range('E28').Value = '  -5.29%  '
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('D30').Value = '''2.25'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.52%  '
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').Value = '''2.72'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -7.72%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').Value = '''0.999'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('D34').Value = '''26.83'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.10%  '
$ws.Range('E35').Value = '  -4.99%  '
$ws.Range('E36').Value = '  -5.87%  '
$ws.Range('D37').Value = '''54.59'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.57%  '
$ws.Range('D38').Value = '''3.25'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.92%  '
$ws.Range('D39').Value = '0.0₃0753'
$ws.Range('E39').Value = '  -3.02%  '
$ws.Range('D40').Value = '''452.01'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -10.09%  '
$ws.Range('E41').Value = '  -5.00%  '
$ws.Range('E42').Value = '  -4.98%  '
$ws.Range('E43').Value = '  -2.50%  '
$ws.Range('D44').Value = '2.879.58'
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = '''2.34'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.05%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').Value = '''0.273'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -8.29%  '
$ws.Range('D47').Value = '''26.75'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -5.88%  '
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('E50').Value = '  -3.62%  '
$ws.Range('D51').Value = '''119.63'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.40%  '
